$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1) onto the new
# header cell (H1) so the new "Save" column header matches the style of
# the other header cells (bold, centered, bordered - style index 1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Set the new header text and the new data value.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
